$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row at position 32 (AZURDUY JUANA 2627) ---
$ws.Rows.Item(32).Insert()
$ws.Range("A32:R32").NumberFormat = "@"
$ws.Range("I32").NumberFormat = "General"
$ws.Range("M32").NumberFormat = "General"
$ws.Range("N32").NumberFormat = "General"
$ws.Range("A32").Value = "7229"
$ws.Range("B32").Value = "9/16/2025"
$ws.Range("C32").Value = "AZURDUY JUANA 2627"
$ws.Range("D32").Value = "13"
$ws.Range("E32").Value = "ICD30814490"
$ws.Range("F32").Value = "NEW"
$ws.Range("G32").Value = "Pendiente de Traspaso PROPIO"
$ws.Range("H32").Value = "Colocar columna para pedir traspaso de nodo propio"
$ws.Range("I32").Value = 1
$ws.Range("J32").Value = "Cambio"
$ws.Range("K32").Value = "Nodo Teco"
$ws.Range("L32").Value = "Pasante"
$ws.Range("M32").Value = -58.469008
$ws.Range("N32").Value = -34.552083
$ws.Range("O32").Value = "Saavedra"
$ws.Range("P32").Value = "Capital Norte"
$ws.Range("Q32").Value = "COG-L"
$ws.Range("R32").Value = "Fuera de Poligono OVL"

# --- Insert new row at position 80 (ARCOS 2263) ---
$ws.Rows.Item(80).Insert()
$ws.Range("A80:R80").NumberFormat = "@"
$ws.Range("I80").NumberFormat = "General"
$ws.Range("M80").NumberFormat = "General"
$ws.Range("N80").NumberFormat = "General"
$ws.Range("A80").Value = "4862"
$ws.Range("B80").Value = "1/23/2025"
$ws.Range("C80").Value = "ARCOS 2263"
$ws.Range("D80").Value = "13"
$ws.Range("E80").Value = "802857379"
$ws.Range("F80").Value = "NEW"
$ws.Range("G80").Value = "Pendiente de Traspaso PROPIO"
$ws.Range("H80").Value = "picada"
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = "Cambio"
$ws.Range("K80").Value = "Nodo Teco"
$ws.Range("L80").Value = "Pasante"
$ws.Range("M80").Value = -58.455082
$ws.Range("N80").Value = -34.558883
$ws.Range("O80").Value = "Saavedra"
$ws.Range("P80").Value = "Capital Norte"
$ws.Range("Q80").Value = "BLO-P"
$ws.Range("R80").Value = "Fuera de Poligono OVL"

# --- Insert new row at position 82 (BARCO CENTENERA DEL 545) ---
$ws.Rows.Item(82).Insert()
$ws.Range("A82:R82").NumberFormat = "@"
$ws.Range("I82").NumberFormat = "General"
$ws.Range("M82").NumberFormat = "General"
$ws.Range("N82").NumberFormat = "General"
$ws.Range("A82").Value = "4528"
$ws.Range("B82").Value = "1/16/2025"
$ws.Range("C82").Value = "BARCO CENTENERA DEL 545"
$ws.Range("D82").Value = "5"
$ws.Range("E82").Value = "802774521"
$ws.Range("F82").Value = "NEW"
$ws.Range("G82").Value = "Pendiente de Traspaso PROPIO"
$ws.Range("H82").Value = "Picada"
$ws.Range("I82").Value = 1
$ws.Range("J82").Value = "Cambio"
$ws.Range("K82").Value = "Sin equipos"
$ws.Range("L82").Value = "Pasante"
$ws.Range("M82").Value = -58.440625
$ws.Range("N82").Value = -34.625499
$ws.Range("O82").Value = "Boedo"
$ws.Range("P82").Value = "Capital Sur"
$ws.Range("Q82").Value = "PCH-C"
$ws.Range("R82").Value = "Fuera de Poligono OVL"

# --- Insert new row at position 86 (DIAZ COLODRERO 3309) ---
$ws.Rows.Item(86).Insert()
$ws.Range("A86:R86").NumberFormat = "@"
$ws.Range("I86").NumberFormat = "General"
$ws.Range("M86").NumberFormat = "General"
$ws.Range("N86").NumberFormat = "General"
$ws.Range("A86").Value = "3299"
$ws.Range("B86").Value = "9/10/2024"
$ws.Range("C86").Value = "DIAZ COLODRERO 3309"
$ws.Range("D86").Value = "12"
$ws.Range("E86").Value = "796186684"
$ws.Range("F86").Value = "NEW"
$ws.Range("G86").Value = "Pendiente de Traspaso PROPIO"
$ws.Range("H86").Value = "qap traspaso nodo TLC y Teco"
$ws.Range("I86").Value = 1
$ws.Range("J86").Value = "Cambio"
$ws.Range("K86").Value = "Nodo Teco"
$ws.Range("L86").Value = "Pasante"
$ws.Range("M86").Value = -58.491722
$ws.Range("N86").Value = -34.565845
$ws.Range("O86").Value = "Paternal"
$ws.Range("P86").Value = "Capital Norte"
$ws.Range("Q86").Value = "PUE-F"
$ws.Range("R86").Value = "Fuera de Poligono OVL"

# --- Insert new row at position 87 (ARCOS 1520) ---
$ws.Rows.Item(87).Insert()
$ws.Range("A87:R87").NumberFormat = "@"
$ws.Range("I87").NumberFormat = "General"
$ws.Range("M87").NumberFormat = "General"
$ws.Range("N87").NumberFormat = "General"
$ws.Range("A87").Value = "5589"
$ws.Range("B87").Value = "12/31/2023"
$ws.Range("C87").Value = "ARCOS 1520"
$ws.Range("D87").Value = "13"
$ws.Range("E87").Value = "799540526"
$ws.Range("F87").Value = "NEW"
$ws.Range("G87").Value = "Pendiente de Traspaso PROPIO"
$ws.Range("H87").Value = "Picada"
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = "Cambio"
$ws.Range("K87").Value = "Nodo Teco"
$ws.Range("L87").Value = "Pasante"
$ws.Range("M87").Value = -58.449125
$ws.Range("N87").Value = -34.565958
$ws.Range("O87").Value = "Colegiales"
$ws.Range("P87").Value = "Capital Norte"
$ws.Range("Q87").Value = "BLO-M"
$ws.Range("R87").Value = "Fuera de Poligono OVL"
